$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text: "second doses" -> "Second Doses" ---
$ws.Range("B1").Value = "Second Doses"

# --- Remove the "Donate" picture/drawing from the sheet ---
if ($ws.Shapes.Count -gt 0) {
    for ($i = $ws.Shapes.Count; $i -ge 1; $i--) {
        $null = $ws.Shapes.Item($i).Delete()
    }
}

# --- Prepend 7 new days of data (16 Nov down to 10 Nov 2021) ---
# Insert 7 blank rows above the current row 2, pushing existing data down.
$null = $ws.Rows("2:8").Insert()

# Copy the number formats/styles from the (now shifted) former row 2 - now row 9 -
# onto the freshly inserted rows so they match the rest of the table exactly.
$ws.Range("A9:B9").Copy()
$ws.Range("A2:B8").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Rows("2:8").RowHeight = 18

$newDates = 44515, 44514, 44513, 44512, 44511, 44510, 44509
$newValues = 4949109, 4940981, 4923992, 4894786, 4865613, 4837627, 4808146

for ($i = 0; $i -lt 7; $i++) {
    $r = 2 + $i
    $ws.Range("A$r").Value = $newDates[$i]
    $ws.Range("B$r").Value = $newValues[$i]
}

# --- Selection moves from E10 to C7 ---
$null = $ws.Range("C7").Select()
